$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.698.18'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '2.059.30'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.664'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.35%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.84'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.40'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.364'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.08%  '
$ws.Range("E11").Value = '  -2.36%  '
$ws.Range("E12").Value = '  -2.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.932'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.64%  '
$ws.Range("E14").Value = '  -3.94%  '
$ws.Range("D15").Value = '2.359.51'
$ws.Range("E16").Value = '  -3.75%  '
$ws.Range("D17").Value = '2.047.43'
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").Value = '36.637.84'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").Value = '0.0₃0865'
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -3.18%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("E31").Value = '  +7.11%  '
$ws.Range("E32").Value = '  -8.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0598'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.38%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0843'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.22'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.35%  '
$ws.Range("E39").Value = '  -3.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.53%  '
$ws.Range("E42").Value = '  -2.68%  '
$ws.Range("E43").Value = '  -3.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0908'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.06%  '
$ws.Range("D46").Value = '1.410.77'
$ws.Range("E46").Value = '  +8.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +14.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.36%  '
$ws.Range("D51").Value = '2.247.99'
$ws.Range("E51").Value = '  +1.15%  '
